$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item('展览')
$wsExpo.Range('F4').Value = 8504
$wsExpo.Range('F5').Value = 8504
$wsExpo.Range('F6').Value = 566
$wsExpo.Range('F7').Value = 7532
$wsExpo.Range('C11').Value = '北京·排球少年同好嘉年华4th'
$wsExpo.Range('D11').Value = '石景山路68号 北京首钢会展中心'
$wsExpo.Range('E11').Value = '2024.07.20 09:00-07.21 17:00'
$wsExpo.Range('F11').Value = 769
$wsExpo.Range('G11').Value = 90
$wsExpo.Range('H11').Value = 'https://show.bilibili.com/platform/detail.html?id=85947'
$wsExpo.Range('I11').Value = '//i2.hdslb.com/bfs/openplatform/202405/bTK0cxNF1716260812082.jpeg'
$wsExpo.Range('B12').NumberFormat = "@"
$wsExpo.Range('B12').Value = '2024-07-20'
$wsExpo.Range('B12').Style = "Normal"
$wsExpo.Range('C12').Value = '北京·英雄的苍穹：正子公也三国、水浒绘画艺术大展'
$wsExpo.Range('D12').Value = '上庄大街18号附近郎园Park（地铁1号线八宝山b口） 郎园Park'
$wsExpo.Range('E12').Value = '2024.07.20 10:00-08.18 19:00'
$wsExpo.Range('F12').Value = 1
$wsExpo.Range('G12').Value = 38
$wsExpo.Range('H12').Value = 'https://show.bilibili.com/platform/detail.html?id=89424'
$wsExpo.Range('I12').Value = '//i2.hdslb.com/bfs/openplatform/202407/hHCrntqE1721180587445.jpeg'
$wsExpo.Range('F14').Value = 192
$wsExpo.Range('F15').Value = 12549
$wsExpo.Range('F18').Value = 2695
$wsExpo.Range('F19').Value = 4893
$wsExpo.Range('F22').Value = 3095
$wsExpo.Range('F24').Value = 140
$wsExpo.Range('F28').Value = 3439
$wsExpo.Range('F29').Value = 82
$wsExpo.Range('F30').Value = 352
$wsExpo.Range('F31').Value = 1788
$wsExpo.Range('F33').Value = 158
$wsExpo.Range('F34').Value = 6228
$wsExpo.Range('F36').Value = 146
$wsExpo.Range('F38').Value = 1891
$wsExpo.Range('F40').Value = 56
$wsExpo.Range('F41').Value = 959
$wsExpo.Range('F42').Value = 7
$wsExpo.Range('F43').Value = 185
$wsExpo.Range('F45').Value = 204
$wsExpo.Range('F46').Value = 1133
$wsExpo.Range('F47').Value = 1122
$wsExpo.Range('F48').Value = 1648
$wsExpo.Range('F50').Value = 128
$wsExpo.Range('G50').Value = 80
$wsShow = $wb.Worksheets.Item('演出')
$wsShow.Range('G6').Value = 98
$wsShow.Range('F14').Value = 115
$wsShow.Range('F20').Value = 83
$wsLocal = $wb.Worksheets.Item('本地生活')
$wsLocal.Range('F2').Value = 368
$wsLocal.Range('F3').Value = 527
$wsAll = $wb.Worksheets.Item('全部类型')
$wsAll.Range('F6').Value = 368
$wsAll.Range('F7').Value = 527
$wsAll.Range('F9').Value = 8504
$wsAll.Range('F10').Value = 566
$wsAll.Range('F11').Value = 7532
$wsAll.Range('F12').Value = 7532
$wsAll.Range('C15').Value = '北京·Summer Overture'
$wsAll.Range('D15').Value = '朝阳北路甲27号菁英梦谷·常营文创产业园南门B5座 WeShow Live 北京'
$wsAll.Range('E15').Value = '2024.07.21 12:00-07.21 19:00'
$wsAll.Range('F15').Value = 261
$wsAll.Range('G15').Value = 98
$wsAll.Range('H15').Value = 'https://show.bilibili.com/platform/detail.html?id=87481'
$wsAll.Range('I15').Value = '//i1.hdslb.com/bfs/openplatform/202406/dP7KKEIk1718608495643.png'
$wsAll.Range('C16').Value = '北京·第三届ICOS X IJOY漫展【七濑公专场见面会】'
$wsAll.Range('D16').Value = '石景山路68号 北京首钢会展中心'
$wsAll.Range('E16').Value = '2024.07.21 12:00-07.21 15:10'
$wsAll.Range('F16').Value = 126
$wsAll.Range('G16').Value = 520
$wsAll.Range('H16').Value = 'https://show.bilibili.com/platform/detail.html?id=87407'
$wsAll.Range('I16').Value = '//i0.hdslb.com/bfs/openplatform/202406/wmNSVTIi1718349284964.jpeg'
$wsAll.Range('F17').Value = 192
$wsAll.Range('F19').Value = 12549
$wsAll.Range('F22').Value = 2695
$wsAll.Range('F23').Value = 2695
$wsAll.Range('F24').Value = 4903
$wsAll.Range('F26').Value = 140
$wsAll.Range('F31').Value = 3439
$wsAll.Range('F32').Value = 352
$wsAll.Range('F33').Value = 1788
$wsAll.Range('F35').Value = 158
$wsAll.Range('F36').Value = 6228
$wsAll.Range('F37').Value = 83
$wsAll.Range('F39').Value = 146
$wsAll.Range('F41').Value = 1891
$wsAll.Range('F44').Value = 56
$wsAll.Range('F45').Value = 959
$wsAll.Range('F46').Value = 185
$wsAll.Range('F47').Value = 204
$wsAll.Range('F48').Value = 1133
$wsAll.Range('F49').Value = 1122
$wsAll.Range('F50').Value = 1648
$wsAll.Range('F52').Value = 128
$wsAll.Range('G52').Value = 80
